$wb = $excel.ActiveWorkbook

function Swap-Rows($ws, $cols) {
    foreach ($col in $cols) {
        $v4 = $ws.Range("$col" + "4").Value2
        $v5 = $ws.Range("$col" + "5").Value2
        $ws.Range("$col" + "4").Value = $v5
        $ws.Range("$col" + "5").Value = $v4
    }
}

$ws1 = $wb.Worksheets.Item(1)
Swap-Rows $ws1 @("A","B","C")

$ws2 = $wb.Worksheets.Item(2)
Swap-Rows $ws2 @("A","B","C","D")

$ws3 = $wb.Worksheets.Item(3)
Swap-Rows $ws3 @("A","B","C","D")

Write-Host "done"
